# Update Fonds de solidarite - 2022-05-13 data refresh
# Applies updated counts/amounts for specific rows (nombre_aides, nombre_entreprises,
# montant_total) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C3").Value = 249326
$ws.Range("E3").Value = 1036473141

$ws.Range("C8").Value = 1045
$ws.Range("D8").Value = 194
$ws.Range("E8").Value = 91109976

$ws.Range("C36").Value = 56967
$ws.Range("E36").Value = 223721899

$ws.Range("C53").Value = 141678
$ws.Range("E53").Value = 590056680

$ws.Range("C54").Value = 60292
$ws.Range("E54").Value = 353611523

$ws.Range("C62").Value = 4188
$ws.Range("E62").Value = 9183365

$ws.Range("C63").Value = 14352
$ws.Range("E63").Value = 36179469

$ws.Range("C64").Value = 5198
$ws.Range("E64").Value = 20337729

$ws.Range("C70").Value = 15724
$ws.Range("E70").Value = 24658131

$ws.Range("C91").Value = 151093
$ws.Range("E91").Value = 482071222

$ws.Range("C92").Value = 408994
$ws.Range("E92").Value = 1593538885

$ws.Range("C93").Value = 209480
$ws.Range("E93").Value = 1307421313

$ws.Range("C94").Value = 94140
$ws.Range("E94").Value = 915436245

$ws.Range("C95").Value = 50709
$ws.Range("E95").Value = 929737465

$ws.Range("C96").Value = 17240
$ws.Range("E96").Value = 789174173

$ws.Range("C98").Value = 809
$ws.Range("E98").Value = 117674774

$ws.Range("C107").Value = 6387
$ws.Range("E107").Value = 21941144

$ws.Range("C109").Value = 1271
$ws.Range("E109").Value = 20757709

$ws.Range("C114").Value = 3791
$ws.Range("E114").Value = 9080369

$ws.Range("C115").Value = 11691
$ws.Range("E115").Value = 32947448

$ws.Range("C116").Value = 4554
$ws.Range("E116").Value = 20454780

$ws.Range("C118").Value = 973
$ws.Range("E118").Value = 11742670

$ws.Range("C122").Value = 8485
$ws.Range("E122").Value = 12669783

$ws.Range("C124").Value = 948
$ws.Range("E124").Value = 2624236

$ws.Range("C142").Value = 168972
$ws.Range("E142").Value = 681760934

$ws.Range("C156").Value = 25104
$ws.Range("E156").Value = 199245925
